$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the rows that changed after repulling data
$ws.Range("F7").Value = 4
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = -6
